$wb = $excel.ActiveWorkbook

# ----- Sheet ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 1092.6
$ws.Range("I28").Value = 1102.8889
$ws.Range("J28").Value = 1000
$ws.Range("K28").Value = 1102.8889
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = -617.8888999999999
$ws.Range("N28").Value = -1970
# Row 51
$ws.Range("H51").Value = 1757.9474
$ws.Range("I51").Value = 1600.0714
$ws.Range("J51").Value = 2200
$ws.Range("K51").Value = 1600.0714
$ws.Range("L51").Value = 2200
$ws.Range("M51").Value = -1116.0714
$ws.Range("N51").Value = -3168
# Row 55
$ws.Range("H55").Value = 531.1
$ws.Range("I55").Value = 320
$ws.Range("J55").Value = 621.5714
$ws.Range("K55").Value = 320
$ws.Range("L55").Value = 621.5714
$ws.Range("M55").Value = -106
$ws.Range("N55").Value = -1049.5714
# Row 62
$ws.Range("H62").Value = 3357.1428
$ws.Range("I62").Value = 1500
$ws.Range("J62").Value = 3666.6667
$ws.Range("K62").Value = 1500
$ws.Range("L62").Value = 3666.6667
$ws.Range("M62").Value = -876
$ws.Range("N62").Value = -4914.6667
# Row 65
$ws.Range("H65").Value = 3357.1428
$ws.Range("I65").Value = 1500
$ws.Range("J65").Value = 3666.6667
$ws.Range("K65").Value = 7500
$ws.Range("L65").Value = 18333.3335
$ws.Range("M65").Value = -4380
$ws.Range("N65").Value = -24573.3335
# Row 74
$ws.Range("H74").Value = 4556.4443
$ws.Range("I74").Value = 3853
$ws.Range("J74").Value = 4908.1665
$ws.Range("K74").Value = 3853
$ws.Range("L74").Value = 4908.1665
$ws.Range("M74").Value = -2917
$ws.Range("N74").Value = -6780.1665
# Row 77
$ws.Range("H77").Value = 4556.4443
$ws.Range("I77").Value = 3853
$ws.Range("J77").Value = 4908.1665
$ws.Range("K77").Value = 19265
$ws.Range("L77").Value = 24540.8325
$ws.Range("M77").Value = -14585
$ws.Range("N77").Value = -33900.8325
# Row 98
$ws.Range("H98").Value = 113911.555
$ws.Range("I98").Value = 202661
$ws.Range("J98").Value = 2974.75
$ws.Range("K98").Value = 202661
$ws.Range("L98").Value = 2974.75
$ws.Range("M98").Value = -201163
$ws.Range("N98").Value = -5970.75
# Row 122
$ws.Range("H122").Value = 113911.555
$ws.Range("I122").Value = 202661
$ws.Range("J122").Value = 2974.75
$ws.Range("K122").Value = 607983
$ws.Range("L122").Value = 8924.25
$ws.Range("M122").Value = -605533
$ws.Range("N122").Value = -13824.25
# Row 135
$ws.Range("H135").Value = 4702.1304
$ws.Range("J135").Value = 590.6667
$ws.Range("L135").Value = 5316.0003
$ws.Range("N135").Value = -10386.0003
# Row 137
$ws.Range("H137").Value = 1904.2222
$ws.Range("J137").Value = 2079.6667
$ws.Range("L137").Value = 6239.000100000001
$ws.Range("N137").Value = -11339.0001

# ----- Sheet ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 9296.877
$ws.Range("I32").Value = 3959.7144
$ws.Range("K32").Value = 3959.7144
$ws.Range("M32").Value = -3672.7144
# Row 122
$ws.Range("H122").Value = 1620.2307
$ws.Range("I122").Value = 1591.04
$ws.Range("J122").Value = 1672.3572
$ws.Range("K122").Value = 4773.12
$ws.Range("L122").Value = 5017.071599999999
$ws.Range("M122").Value = -2323.12
$ws.Range("N122").Value = -9917.071599999999

# ----- Sheet BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 1347.1428
$ws.Range("I107").Value = 1275.8823
$ws.Range("J107").Value = 1650
$ws.Range("K107").Value = 1275.8823
$ws.Range("L107").Value = 1650
$ws.Range("M107").Value = 644.1177
$ws.Range("N107").Value = -5490
# Row 110
$ws.Range("H110").Value = 52500
$ws.Range("J110").Value = 52500
$ws.Range("L110").Value = 52500
$ws.Range("N110").Value = -60680

# ----- Sheet CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1529.8
$ws.Range("I31").Value = 847.9666999999999
$ws.Range("K31").Value = 847.9666999999999
$ws.Range("M31").Value = -552.9666999999999
# Row 34
$ws.Range("H34").Value = 1529.8
$ws.Range("I34").Value = 847.9666999999999
$ws.Range("K34").Value = 847.9666999999999
$ws.Range("M34").Value = -645.9666999999999
# Row 50
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
# Row 86
$ws.Range("H86").Value = 163819.27
$ws.Range("I86").Value = 276586.16
$ws.Range("J86").Value = 2723.7144
$ws.Range("K86").Value = 276586.16
$ws.Range("L86").Value = 2723.7144
$ws.Range("M86").Value = -275463.16
$ws.Range("N86").Value = -4969.7144
# Row 89
$ws.Range("H89").Value = 163819.27
$ws.Range("I89").Value = 276586.16
$ws.Range("J89").Value = 2723.7144
$ws.Range("K89").Value = 1382930.8
$ws.Range("L89").Value = 13618.572
$ws.Range("M89").Value = -1377314.8
$ws.Range("N89").Value = -24850.572
# Row 112
$ws.Range("H112").Value = 67000
$ws.Range("J112").Value = 67000
$ws.Range("L112").Value = 67000
$ws.Range("N112").Value = -69954

# ----- Sheet CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 98
$ws.Range("H98").Value = 367.8
$ws.Range("J98").Value = 359.75
$ws.Range("L98").Value = 1079.25
$ws.Range("N98").Value = -4075.25

# ----- Sheet GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 6
$ws.Range("H6").Value = 31666.666
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 31666.666
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 31666.666
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -31892.666
# Row 16
$ws.Range("H16").Value = 31666.666
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 31666.666
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 31666.666
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -32166.666
# Row 28
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
# Row 122
$ws.Range("H122").Value = 1924.7826
$ws.Range("I122").Value = 1817.4
$ws.Range("J122").Value = 2007.3846
$ws.Range("K122").Value = 5452.200000000001
$ws.Range("L122").Value = 6022.1538
$ws.Range("M122").Value = -3002.200000000001
$ws.Range("N122").Value = -10922.1538
# Row 132
$ws.Range("H132").Value = 7366.304
$ws.Range("I132").Value = 8525.8125
$ws.Range("J132").Value = 4716
$ws.Range("K132").Value = 25577.4375
$ws.Range("L132").Value = 14148
$ws.Range("M132").Value = -23047.4375
$ws.Range("N132").Value = -19208

# ----- Sheet LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2475.4285
$ws.Range("I7").Value = 2052.2666
$ws.Range("J7").Value = 3533.3333
$ws.Range("K7").Value = 2052.2666
$ws.Range("L7").Value = 3533.3333
$ws.Range("M7").Value = -1940.2666
$ws.Range("N7").Value = -3757.3333
# Row 22
$ws.Range("H22").Value = 709.1111
$ws.Range("I22").Value = 668.5714
$ws.Range("J22").Value = 851
$ws.Range("K22").Value = 668.5714
$ws.Range("L22").Value = 851
$ws.Range("M22").Value = -373.5714
$ws.Range("N22").Value = -1441
# Row 27
$ws.Range("H27").Value = 709.1111
$ws.Range("I27").Value = 668.5714
$ws.Range("J27").Value = 851
$ws.Range("K27").Value = 668.5714
$ws.Range("L27").Value = 851
$ws.Range("M27").Value = -561.5714
$ws.Range("N27").Value = -1065
# Row 40
$ws.Range("H40").Value = 2822.9429
$ws.Range("I40").Value = 2582.0356
$ws.Range("K40").Value = 2582.0356
$ws.Range("M40").Value = -2446.0356
# Row 46
$ws.Range("H46").Value = 1541.8334
$ws.Range("I46").Value = 1337.75
$ws.Range("K46").Value = 1337.75
$ws.Range("M46").Value = -1149.75
# Row 93
$ws.Range("H93").Value = 1217.4
$ws.Range("I93").Value = 1029.04
$ws.Range("J93").Value = 1531.3334
$ws.Range("K93").Value = 1029.04
$ws.Range("L93").Value = 1531.3334
$ws.Range("M93").Value = 218.96
$ws.Range("N93").Value = -4027.3334
# Row 126
$ws.Range("H126").Value = 2475.4285
$ws.Range("I126").Value = 2052.2666
$ws.Range("J126").Value = 3533.3333
$ws.Range("K126").Value = 6156.7998
$ws.Range("L126").Value = 10599.9999
$ws.Range("M126").Value = -3686.7998
$ws.Range("N126").Value = -15539.9999

# ----- Sheet WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 29
$ws.Range("H29").Value = 3015
$ws.Range("J29").Value = 3015
$ws.Range("L29").Value = 3015
$ws.Range("N29").Value = -3595
# Row 122
$ws.Range("H122").Value = 1849.8636
$ws.Range("I122").Value = 1412.3334
$ws.Range("J122").Value = 3818.75
$ws.Range("K122").Value = 4237.0002
$ws.Range("L122").Value = 11456.25
$ws.Range("M122").Value = -1787.0002
$ws.Range("N122").Value = -16356.25
# Row 132
$ws.Range("H132").Value = 1582.0344
$ws.Range("I132").Value = 1110.0555
$ws.Range("J132").Value = 2354.3635
$ws.Range("K132").Value = 3330.1665
$ws.Range("L132").Value = 7063.0905
$ws.Range("M132").Value = -800.1664999999998
$ws.Range("N132").Value = -12123.0905

Write-Host "Applied market-data refresh across 8 sheets"

